$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 150.9
$ws.Cells.Item(18, 9).Value = 150.9
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 150.9
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 133.1

$ws.Cells.Item(126, 8).Value = 11989
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 11989
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 11989
$ws.Cells.Item(126, 14).Value = -21869

$ws.Cells.Item(129, 8).Value = 847.38336
$ws.Cells.Item(129, 9).Value = 597.6
$ws.Cells.Item(129, 10).Value = 870.0909
$ws.Cells.Item(129, 11).Value = 1792.8
$ws.Cells.Item(129, 12).Value = 2610.2727
$ws.Cells.Item(129, 13).Value = 3207.2
$ws.Cells.Item(129, 14).Value = -12610.2727

$ws.Cells.Item(138, 8).Value = 2845.6904
$ws.Cells.Item(138, 9).Value = 776.8
$ws.Cells.Item(138, 10).Value = 3125.2703
$ws.Cells.Item(138, 11).Value = 2330.4
$ws.Cells.Item(138, 12).Value = 9375.8109
$ws.Cells.Item(138, 13).Value = 2809.6
$ws.Cells.Item(138, 14).Value = -19655.8109

$ws.Cells.Item(141, 8).Value = 1615.8334
$ws.Cells.Item(141, 9).Value = 942.8125
$ws.Cells.Item(141, 10).Value = 7000
$ws.Cells.Item(141, 11).Value = 2828.4375
$ws.Cells.Item(141, 12).Value = 21000
$ws.Cells.Item(141, 13).Value = 2351.5625
$ws.Cells.Item(141, 14).Value = -31360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(7, 8).Value = 50000
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 50000
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 50000
$ws.Cells.Item(7, 14).Value = -50228

$ws.Cells.Item(110, 8).Value = 1570.5217
$ws.Cells.Item(110, 9).Value = 1302.625
$ws.Cells.Item(110, 10).Value = 2182.8572
$ws.Cells.Item(110, 11).Value = 1302.625
$ws.Cells.Item(110, 12).Value = 2182.8572
$ws.Cells.Item(110, 13).Value = 742.375
$ws.Cells.Item(110, 14).Value = -6272.8572

$ws.Cells.Item(132, 8).Value = 58924.668
$ws.Cells.Item(132, 9).Value = 3256
$ws.Cells.Item(132, 10).Value = 103459.6
$ws.Cells.Item(132, 11).Value = 9768
$ws.Cells.Item(132, 12).Value = 310378.8
$ws.Cells.Item(132, 13).Value = -7238
$ws.Cells.Item(132, 14).Value = -315438.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 547.8333
$ws.Cells.Item(107, 9).Value = 572.88464
$ws.Cells.Item(107, 10).Value = 385
$ws.Cells.Item(107, 11).Value = 572.88464
$ws.Cells.Item(107, 12).Value = 385
$ws.Cells.Item(107, 13).Value = 1347.11536
$ws.Cells.Item(107, 14).Value = -4225

$ws.Cells.Item(134, 8).Value = 22313.611
$ws.Cells.Item(134, 9).Value = 24691.523
$ws.Cells.Item(134, 10).Value = 1388
$ws.Cells.Item(134, 11).Value = 74074.569
$ws.Cells.Item(134, 12).Value = 4164
$ws.Cells.Item(134, 13).Value = -71539.569
$ws.Cells.Item(134, 14).Value = -9234

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 17000
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 17000
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 17000
$ws.Cells.Item(50, 14).Value = -18250

$ws.Cells.Item(58, 8).Value = 18268
$ws.Cells.Item(58, 9).Value = 1145.6364
$ws.Cells.Item(58, 10).Value = 72081.14
$ws.Cells.Item(58, 11).Value = 1145.6364
$ws.Cells.Item(58, 12).Value = 72081.14
$ws.Cells.Item(58, 13).Value = -942.6364000000001
$ws.Cells.Item(58, 14).Value = -72487.14

$ws.Cells.Item(59, 8).Value = 21200
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 21200
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 12).Value = 21200
$ws.Cells.Item(59, 14).Value = -23490

$ws.Cells.Item(60, 8).Value = 15000
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 15000
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 12).Value = 15000
$ws.Cells.Item(60, 14).Value = -16022

$ws.Cells.Item(62, 8).Value = 4628.5
$ws.Cells.Item(62, 9).Value = 3003.3333
$ws.Cells.Item(62, 10).Value = 5603.6
$ws.Cells.Item(62, 11).Value = 3003.3333
$ws.Cells.Item(62, 12).Value = 5603.6
$ws.Cells.Item(62, 13).Value = -2379.3333
$ws.Cells.Item(62, 14).Value = -6851.6

$ws.Cells.Item(65, 8).Value = 4628.5
$ws.Cells.Item(65, 9).Value = 3003.3333
$ws.Cells.Item(65, 10).Value = 5603.6
$ws.Cells.Item(65, 11).Value = 15016.6665
$ws.Cells.Item(65, 12).Value = 28018
$ws.Cells.Item(65, 13).Value = -11896.6665
$ws.Cells.Item(65, 14).Value = -34258

$ws.Cells.Item(132, 8).Value = 20589.172
$ws.Cells.Item(132, 9).Value = 28324.947
$ws.Cells.Item(132, 10).Value = 5891.2
$ws.Cells.Item(132, 11).Value = 84974.841
$ws.Cells.Item(132, 12).Value = 17673.6
$ws.Cells.Item(132, 13).Value = -82444.841
$ws.Cells.Item(132, 14).Value = -22733.6

$ws.Cells.Item(136, 8).Value = 18268
$ws.Cells.Item(136, 9).Value = 1145.6364
$ws.Cells.Item(136, 10).Value = 72081.14
$ws.Cells.Item(136, 11).Value = 3436.9092
$ws.Cells.Item(136, 12).Value = 216243.42
$ws.Cells.Item(136, 13).Value = -886.9092000000001
$ws.Cells.Item(136, 14).Value = -221343.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1285.625
$ws.Cells.Item(68, 9).Value = 674.5
$ws.Cells.Item(68, 10).Value = 1341.1818
$ws.Cells.Item(68, 11).Value = 2023.5
$ws.Cells.Item(68, 12).Value = 4023.5454
$ws.Cells.Item(68, 13).Value = -1212.5
$ws.Cells.Item(68, 14).Value = -5645.5454

$ws.Cells.Item(71, 8).Value = 1285.625
$ws.Cells.Item(71, 9).Value = 674.5
$ws.Cells.Item(71, 10).Value = 1341.1818
$ws.Cells.Item(71, 11).Value = 6070.5
$ws.Cells.Item(71, 12).Value = 12070.6362
$ws.Cells.Item(71, 13).Value = -2014.5
$ws.Cells.Item(71, 14).Value = -20182.6362

$ws.Cells.Item(81, 8).Value = 4006.9167
$ws.Cells.Item(81, 9).Value = 756.5
$ws.Cells.Item(81, 10).Value = 4657
$ws.Cells.Item(81, 11).Value = 2269.5
$ws.Cells.Item(81, 12).Value = 13971
$ws.Cells.Item(81, 13).Value = -1146.5
$ws.Cells.Item(81, 14).Value = -16217

$ws.Cells.Item(84, 8).Value = 4006.9167
$ws.Cells.Item(84, 9).Value = 756.5
$ws.Cells.Item(84, 10).Value = 4657
$ws.Cells.Item(84, 11).Value = 6808.5
$ws.Cells.Item(84, 12).Value = 41913
$ws.Cells.Item(84, 13).Value = -1192.5
$ws.Cells.Item(84, 14).Value = -53145

$ws.Cells.Item(107, 8).Value = 4292.8887
$ws.Cells.Item(107, 9).Value = 8819.666999999999
$ws.Cells.Item(107, 10).Value = 671.4666999999999
$ws.Cells.Item(107, 11).Value = 26459.001
$ws.Cells.Item(107, 12).Value = 2014.4001
$ws.Cells.Item(107, 13).Value = -24539.001
$ws.Cells.Item(107, 14).Value = -5854.4001

$ws.Cells.Item(131, 8).Value = 164781.03
$ws.Cells.Item(131, 9).Value = 1015
$ws.Cells.Item(131, 10).Value = 176273.39
$ws.Cells.Item(131, 11).Value = 3045
$ws.Cells.Item(131, 12).Value = 528820.17
$ws.Cells.Item(131, 13).Value = 1995
$ws.Cells.Item(131, 14).Value = -538900.17

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4354.7
$ws.Cells.Item(70, 9).Value = 3959.6
$ws.Cells.Item(70, 10).Value = 4749.8
$ws.Cells.Item(70, 11).Value = 3959.6
$ws.Cells.Item(70, 12).Value = 4749.8
$ws.Cells.Item(70, 13).Value = -3689.6
$ws.Cells.Item(70, 14).Value = -5289.8

$ws.Cells.Item(73, 8).Value = 4354.7
$ws.Cells.Item(73, 9).Value = 3959.6
$ws.Cells.Item(73, 10).Value = 4749.8
$ws.Cells.Item(73, 11).Value = 3959.6
$ws.Cells.Item(73, 12).Value = 4749.8
$ws.Cells.Item(73, 13).Value = -3023.6
$ws.Cells.Item(73, 14).Value = -6621.8

$ws.Cells.Item(96, 8).Value = 19261
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 19261
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 12).Value = 19261
$ws.Cells.Item(96, 14).Value = -24753

$ws.Cells.Item(122, 8).Value = 3573.15
$ws.Cells.Item(122, 9).Value = 2482.5715
$ws.Cells.Item(122, 10).Value = 6117.8335
$ws.Cells.Item(122, 11).Value = 7447.7145
$ws.Cells.Item(122, 12).Value = 18353.5005
$ws.Cells.Item(122, 13).Value = -4997.7145
$ws.Cells.Item(122, 14).Value = -23253.5005

$ws.Cells.Item(132, 8).Value = 45806.03
$ws.Cells.Item(132, 9).Value = 38315.074
$ws.Cells.Item(132, 10).Value = 74699.71000000001
$ws.Cells.Item(132, 11).Value = 114945.222
$ws.Cells.Item(132, 12).Value = 224099.13
$ws.Cells.Item(132, 13).Value = -112415.222
$ws.Cells.Item(132, 14).Value = -229159.13

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3669.6667
$ws.Cells.Item(7, 9).Value = 3961.7646
$ws.Cells.Item(7, 10).Value = 2960.2856
$ws.Cells.Item(7, 11).Value = 3961.7646
$ws.Cells.Item(7, 12).Value = 2960.2856
$ws.Cells.Item(7, 13).Value = -3849.7646
$ws.Cells.Item(7, 14).Value = -3184.2856

$ws.Cells.Item(16, 8).Value = 957.8182
$ws.Cells.Item(16, 9).Value = 732.2222
$ws.Cells.Item(16, 10).Value = 1973
$ws.Cells.Item(16, 11).Value = 732.2222
$ws.Cells.Item(16, 12).Value = 1973
$ws.Cells.Item(16, 13).Value = -562.2222
$ws.Cells.Item(16, 14).Value = -2313

$ws.Cells.Item(46, 8).Value = 1304.0476
$ws.Cells.Item(46, 9).Value = 944.5833
$ws.Cells.Item(46, 10).Value = 1783.3334
$ws.Cells.Item(46, 11).Value = 944.5833
$ws.Cells.Item(46, 12).Value = 1783.3334
$ws.Cells.Item(46, 13).Value = -756.5833
$ws.Cells.Item(46, 14).Value = -2159.3334

$ws.Cells.Item(68, 8).Value = 7672.727
$ws.Cells.Item(68, 9).Value = 3450
$ws.Cells.Item(68, 10).Value = 8611.111000000001
$ws.Cells.Item(68, 11).Value = 3450
$ws.Cells.Item(68, 12).Value = 8611.111000000001
$ws.Cells.Item(68, 13).Value = -2701
$ws.Cells.Item(68, 14).Value = -10109.111

$ws.Cells.Item(71, 8).Value = 7672.727
$ws.Cells.Item(71, 9).Value = 3450
$ws.Cells.Item(71, 10).Value = 8611.111000000001
$ws.Cells.Item(71, 11).Value = 17250
$ws.Cells.Item(71, 12).Value = 43055.55500000001
$ws.Cells.Item(71, 13).Value = -13506
$ws.Cells.Item(71, 14).Value = -50543.55500000001

$ws.Cells.Item(126, 8).Value = 3669.6667
$ws.Cells.Item(126, 9).Value = 3961.7646
$ws.Cells.Item(126, 10).Value = 2960.2856
$ws.Cells.Item(126, 11).Value = 11885.2938
$ws.Cells.Item(126, 12).Value = 8880.856800000001
$ws.Cells.Item(126, 13).Value = -9415.293799999999
$ws.Cells.Item(126, 14).Value = -13820.8568

$ws.Cells.Item(135, 8).Value = 23214.5
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 23214.5
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 23214.5
$ws.Cells.Item(135, 14).Value = -33354.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1165.9056
$ws.Cells.Item(132, 9).Value = 888.75
$ws.Cells.Item(132, 10).Value = 2018.6923
$ws.Cells.Item(132, 11).Value = 2666.25
$ws.Cells.Item(132, 12).Value = 6056.0769
$ws.Cells.Item(132, 13).Value = -136.25
$ws.Cells.Item(132, 14).Value = -11116.0769
